{"js": "// Locate the paragraph that contains the \"History (list of steps done)...\" bullet\n// and wrap its text with \"[ \" ... \" ]  -  isn't it the same as \"Work done\"?\"\n// while keeping the original sentence as its own run (matching the target OOXML\n// diff, which splits the text into three separate runs within the same\n// paragraph).\n\nconst targetText = \"History (list of steps done), user can undo any step\";\n\nconst results = context.document.body.search(targetText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target paragraph text: \" + targetText);\n}\n\nconst hit = results.items[0];\n\n// Insert the leading bracket right before the matched text.\nhit.insertText(\"[ \", \"Before\");\n\n// Insert the trailing annotation right after the matched text.\nhit.insertText(\" ]  -  isn\\u2019t it the same as \\u201CWork done\\u201D?\", \"After\");\n\nawait context.sync();\n", "ps1": "# Wrap the \"History (list of steps done), user can undo any step\" bullet\n# with \"[ \" ... \" ]  -  isn't it the same as \"Work done\"?\" as in the target\n# revision (the sentence itself is left untouched, text is only added\n# immediately before and immediately after it).\n\n$d = $word.ActiveDocument\n\n$targetText = \"History (list of steps done), user can undo any step\"\n$suffixText = \" ]  -  isn\" + [char]0x2019 + \"t it the same as \" + [char]0x201C + \"Work done\" + [char]0x201D + \"?\"\n\n# Append the trailing annotation right after the sentence.\n$rngEnd = $d.Content\n$rngEnd.Find.ClearFormatting()\n$rngEnd.Find.Execute($targetText) | Out-Null\n$rngEnd.Collapse(0)  # wdCollapseEnd\n$rngEnd.InsertAfter($suffixText)\n\n# Prepend the leading bracket right before the sentence.\n$rngStart = $d.Content\n$rngStart.Find.ClearFormatting()\n$rngStart.Find.Execute($targetText) | Out-Null\n$rngStart.Collapse(1)  # wdCollapseStart\n$rngStart.InsertBefore(\"[ \")\n"}
